# personnel edits to fix roles, spacing, acknowledgments
$wb = $excel.ActiveWorkbook

$personnel = $wb.Worksheets.Item("Personnel")

# Fix role spelling/spacing for Rachel Stanley (principal investigator)
$personnel.Range("G5").Value = "principal Investigator"

# Fix role spelling/spacing for Jaxine Wolfe and Kate Morkeski (metadata providers)
$personnel.Range("G7").Value = "metadata Provider"
$personnel.Range("G8").Value = "metadata Provider"

# Update active sheet / selection to match the author's final view state
$personnel.Activate()
$personnel.Range("C13").Select()
